$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values that must stay stored as TEXT
# (the sheet uses inline/shared strings throughout). Prefixing with a
# leading apostrophe forces Excel to keep the entry as text instead of
# auto-converting it to a Number.
$ws.Range("D2").Value  = "'244.35"
$ws.Range("D3").Value  = "'23.97"
$ws.Range("D4").Value  = "'5.261"
$ws.Range("D5").Value  = "'0.05831"
$ws.Range("D6").Value  = "'6.458"
$ws.Range("D7").Value  = "'3.230"
$ws.Range("D8").Value  = "'0.8081"
$ws.Range("D9").Value  = "'0.8871"
$ws.Range("D11").Value = "'0.07107"
$ws.Range("D12").Value = "'0.03087"
$ws.Range("D14").Value = "'0.09336"
$ws.Range("D15").Value = "'3.830"
$ws.Range("D16").Value = "'0.001541"
$ws.Range("D17").Value = "'0.04707"
$ws.Range("D18").Value = "'0.0006051"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006148"
$ws.Range("D20").Value = "'0.001257"
$ws.Range("D22").Value = "'0.00008700"
$ws.Range("D23").Value = "'3.544"
$ws.Range("D40").Value = "'0.03838"

# Rows 41-43 got reshuffled (BKEXToken / CEJI / KickToken rotated)
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006280"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1052"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002541"
$ws.Range("E43").Value = "42CEJICEJI"

$ws.Range("D44").Value = "'0.007236"
$ws.Range("D45").Value = "'0.00005338"
$ws.Range("D47").Value = "'0.5214"
$ws.Range("D48").Value = "'0.003389"
